$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Type"
$ws.Range("B1").Value = "Weight"
$ws.Range("C1").Value = "Group"

# Row 2 - SN1 fatty acid
$ws.Range("A2").Value = "SN1_[FA-H]-"
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 1

# Row 3 - SN2 fatty acid
$ws.Range("A3").Value = "SN2_[FA-H]-"
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 1

# Row 4 - LPL(SN1) minus H
$ws.Range("A4").Value = "[LPL(SN1)-H]-"
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 2

# Row 5 - LPL(SN2) minus H
$ws.Range("A5").Value = "[LPL(SN2)-H]-"
$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 2

# Row 6 - LPL(SN1) minus H2O minus H
$ws.Range("A6").Value = "[LPL(SN1)-H2O-H]-"
$ws.Range("B6").Value = 10
$ws.Range("C6").Value = 2

# Row 7 - LPL(SN2) minus H2O minus H
$ws.Range("A7").Value = "[LPL(SN2)-H2O-H]-"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 2

# Drop the now-unused D:F columns (table shrinks from 6 to 3 columns)
$ws.Range("D1:F7").EntireColumn.Delete()

# Match the author's final cursor position
[void]$ws.Range("C3").Select()
